$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was recorded for "Jengibre" at
# "Terminal La Palmera de La Serena". It belongs right after the existing
# row 82, so push the existing row 83 (and everything below it) down by
# one row, then fill the freshly opened row 83 with the new record.
$ws.Rows.Item(83).Insert()

$ws.Cells.Item(83, 1).Value = 8
$ws.Cells.Item(83, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = 45068
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = 100114007
$ws.Cells.Item(83, 7).Value = "Jengibre"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 460
$ws.Cells.Item(83, 11).Value = 17000
$ws.Cells.Item(83, 12).Value = 18000
$ws.Cells.Item(83, 13).Value = 17500
$ws.Cells.Item(83, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(83, 15).Value = "Perú"
$ws.Cells.Item(83, 16).Value = 1346
$ws.Cells.Item(83, 17).Value = 13
$ws.Cells.Item(83, 18).Value = "Hortaliza"
